$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.081.15"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "2.104.11"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.71%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "345.57"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("E6").Value = "  -0.62%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.5179"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -1.70%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.4433"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +4.23%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "52.39"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -2.60%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "1.176"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "25.29"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +3.56%  "
$ws.Range("D13").Value = "2.108.24"
$ws.Range("E13").Value = "  -0.76%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "6.722"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "8.121"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "99.65"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +1.26%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.00001170"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "20.73"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +6.26%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.06703"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -0.63%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "6.223"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").Value = "30.163.73"
$ws.Range("E23").Value = "  -1.86%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "12.68"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "2.332"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").Value = "2.357.87"
$ws.Range("E26").Value = "  -0.47%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "22.04"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -1.72%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "164.75"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "2.548"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "133.68"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.166"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("E32").Value = "  -1.62%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.636"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "6.251"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -2.34%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "3.966"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "6.211"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +4.51%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "10.15"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -3.48%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.02569"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -3.81%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.06789"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.2285"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.6949"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +0.85%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "12.54"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -0.65%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "1.308"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +3.46%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.6680"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +3.36%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "14.23"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -6.54%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "2.282"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.16%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "3.643"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -1.57%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.00000000357"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -2.59%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "1.222"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "82.51"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.07196"
$r.Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "
